$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
$ws.Range("D14").ClearContents()
